# Applies updated model metrics (C/D/E columns, rows 2-9) and refreshes the
# per-column Greens background-gradient shading on the RMSE (D) and U (E) columns.
# NOTE: Range.Interior.Color takes a COLORREF (0xBBGGRR), so hex literals below
# are the byte-swapped form of the target RRGGBB fill color.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = -8.163399999999999
$ws.Range("D2").Value = 0.8256
$ws.Range("D2").Interior.Color = 0xF5FCF7  # fill F7FCF5
$ws.Range("E2").Value = 2.629
$ws.Range("E2").Interior.Color = 0xF5FCF7  # fill F7FCF5

# Row 3
$ws.Range("C3").Value = -1.5348
$ws.Range("D3").Value = 0.6611
$ws.Range("D3").Interior.Color = 0xB4E4BB  # fill BBE4B4
$ws.Range("E3").Value = 1.5656
$ws.Range("E3").Interior.Color = 0x70BE68  # fill 68BE70

# Row 4
$ws.Range("C4").Value = -0.0432
$ws.Range("D4").Value = 0.5663
$ws.Range("D4").Interior.Color = 0x83CC84  # fill 84CC83
$ws.Range("E4").Value = 1.3511
$ws.Range("E4").Interior.Color = 0x5AA73E  # fill 3EA75A

# Row 5
$ws.Range("C5").Value = 0.8611
$ws.Range("D5").Value = 0.2582
$ws.Range("D5").Interior.Color = 0x1B4400  # fill 00441B
$ws.Range("E5").Value = 0.6246
$ws.Range("E5").Interior.Color = 0x1B4400  # fill 00441B

# Row 6
$ws.Range("C6").Value = 0.6677999999999999
$ws.Range("D6").Value = 0.4465
$ws.Range("D6").Interior.Color = 0x549F36  # fill 369F54
$ws.Range("E6").Value = 1.1071
$ws.Range("E6").Interior.Color = 0x438820  # fill 208843

# Row 7
$ws.Range("C7").Value = 0.5405
$ws.Range("D7").Value = 0.5273
$ws.Range("D7").Interior.Color = 0x71BF6A  # fill 6ABF71
$ws.Range("E7").Value = 1.3321
$ws.Range("E7").Interior.Color = 0x59A53C  # fill 3CA559

# Row 8
$ws.Range("C8").Value = 0.2033
$ws.Range("D8").Value = 0.7059
$ws.Range("D8").Interior.Color = 0xCAEDD0  # fill D0EDCA
$ws.Range("E8").Value = 1.7616
$ws.Range("E8").Interior.Color = 0x8AD08D  # fill 8DD08A

# Row 9
$ws.Range("C9").Value = 0.0177
$ws.Range("D9").Value = 0.7993
$ws.Range("D9").Interior.Color = 0xEEFAF1  # fill F1FAEE
$ws.Range("E9").Value = 1.9911
$ws.Range("E9").Interior.Color = 0xACE0B2  # fill B2E0AC
